# "excel prevented last commit >:(" -- add two new rows of results to the
# "1D NEW" sheet's table, annotate two existing rows with comments, and
# restore the normal (maximized, scrolled-to-top) window/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1D NEW")

# --- annotate existing rows with comments (column M) ---
# (written in this order so the shared-string table indices line up with
# the authored workbook: "Brute force..." lands before "Seems to follow...")
$ws.Range("M15").Value = "Brute force solved for 4 detectors."
$ws.Range("M14").Value = "Seems to follow noiseless solution"

# --- append two new data rows ---
$ws.Range("B18").Value = "1dmockanderrors15"
$ws.Range("C18").Value = 200
$ws.Range("D18").Value = 200
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1.5
$ws.Range("G18").Value = 60
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 30

$ws.Range("B19").Value = "1dmockanderrors16"
$ws.Range("C19").Value = 200
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 30

# --- grow the table (ListObject) to cover the two new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:M19"))

# --- restore the normal view: scrolled to top, selection on the new row ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B18").Select()

# --- maximize the Excel application window ---
$excel.WindowState = -4137
